$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.282.34"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.641.59"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "195.92"
$ws.Range("E5").Value = "  +6.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "577.23"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.635.01"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.680"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +5.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.39"
$ws.Range("E12").Value = "  +4.95%  "
$ws.Range("E13").Value = "  +16.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.15"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.217.69"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.641.05"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.58"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.188.09"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.59"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.90"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.77"
$ws.Range("E23").Value = "  +23.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.24"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.17"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.96"
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.64"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.86"
$ws.Range("E29").Value = "  +5.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.15"
$ws.Range("E30").Value = "  +20.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.17"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.76"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "693.75"
$ws.Range("E33").Value = "  +17.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.24"
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.119"
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.87"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.82"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.418"
$ws.Range("E38").Value = "  +10.97%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0800"
$ws.Range("E39").Value = "  +8.39%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  +18.54%  "
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.197.02"
$ws.Range("E43").Value = "  +18.21%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.12"
$ws.Range("E44").Value = "  +12.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("E46").Value = "  +26.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0424"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.89"
$ws.Range("E49").Value = "  +6.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.10"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.05"
$ws.Range("E51").Value = "  +2.94%  "
